$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 86

$ws.Cells.Item($newRow, 1).Value = 7
$ws.Cells.Item($newRow, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($newRow, 3).Value = "Ñuble"
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 4).Value = 44890
$ws.Cells.Item($newRow, 5).Value = 16
$ws.Cells.Item($newRow, 6).Value = 100112022
$ws.Cells.Item($newRow, 7).Value = "Arveja Verde"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 80
$ws.Cells.Item($newRow, 11).Value = 19000
$ws.Cells.Item($newRow, 12).Value = 20000
$ws.Cells.Item($newRow, 13).Value = 19500
$ws.Cells.Item($newRow, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item($newRow, 15).Value = "Región del Maule"
$ws.Cells.Item($newRow, 16).Value = 780
$ws.Cells.Item($newRow, 17).Value = 25
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
